$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction")

# --- Rebuild the Introduction sheet content row by row ---

# Row 1: Title "Help regarding this Configuration File" - taller row, wrapped text
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Cells.Item(1,1).Value2 = "Help regarding this Configuration File"
$ws.Cells.Item(1,1).WrapText = $true

# Row 2: blank spacer row (wrap text style, no content)
$ws.Cells.Item(2,1).ClearContents()
$ws.Cells.Item(2,1).WrapText = $true

# Row 3: Bold legend heading "####  Legend of Key Value pairs####"
$ws.Cells.Item(3,1).Value2 = "####  Legend of Key Value pairs####"
$ws.Cells.Item(3,1).Font.Bold = $true
$ws.Cells.Item(3,1).Font.Size = 14
$ws.Cells.Item(3,1).WrapText = $true

# Row 4: explanatory paragraph, taller row
$ws.Rows.Item(4).RowHeight = 30
$ws.Cells.Item(4,1).Value2 = "You may want to mark keys in this settings dictionary with certain colours. One use I've needed for colors was to define the keys needed in the operation of the framework. "
$ws.Cells.Item(4,1).WrapText = $true

# Row 5: Input-style legend entry
$ws.Cells.Item(5,1).Value2 = "This key is used in the Framework layer. You can change the values, but do not delete the keys"
$ws.Cells.Item(5,1).Style = "Input"
$ws.Cells.Item(5,1).WrapText = $true

# Row 6: Good-style legend entry, taller row
$ws.Rows.Item(6).RowHeight = 30
$ws.Cells.Item(6,1).Value2 = "This key is used in the Business Process Layer. The developer is responsible for the keys. The user is responsible for the values."
$ws.Cells.Item(6,1).Style = "Good"
$ws.Cells.Item(6,1).WrapText = $true

# Row 7: Note-style legend entry (bold)
$ws.Cells.Item(7,1).Value2 = "This key belongs to user designated category 1"
$ws.Cells.Item(7,1).Style = "Note"
$ws.Cells.Item(7,1).Font.Bold = $true
$ws.Cells.Item(7,1).WrapText = $true

# Row 8: Bad-style legend entry
$ws.Cells.Item(8,1).Value2 = "This Key is used only during Debug. You can delete in production"
$ws.Cells.Item(8,1).Style = "Bad"
$ws.Cells.Item(8,1).WrapText = $true

# Row 10 previously had an explicit taller height (ht=30) - restore default
$ws.Rows.Item(10).AutoFit()

# Rows 9-14: blank trailing rows
for ($r = 9; $r -le 14; $r++) {
  $ws.Cells.Item($r,1).ClearContents()
  $ws.Cells.Item($r,1).Style = "Normal"
  $ws.Cells.Item($r,1).WrapText = $false
}

# Selection moves to A10
$ws.Range("A10").Select()
